$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new region/location rows after the existing data (rows 1-8 already
# contain the header + 7 entries). Add two more entries so the DFs reflect
# the latest stock locations.
$ws.Range("A9").Value = "天上"
$ws.Range("B9").Value = "白云"
$ws.Range("A10").Value = "我的世界"
$ws.Range("B10").Value = "下界"

$wb.Save()
